$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Cells.Item(6, 1).Value = "F35363_1556815384"
$ws.Cells.Item(6, 2).Value = 43587.73827284455
$ws.Cells.Item(6, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(6, 3).Value = "'123"
$ws.Cells.Item(6, 4).Value = "Carteiras"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = "'01"
$ws.Cells.Item(6, 7).Value = "Marca: Herschel; Tamanho: Unic; "
$ws.Cells.Item(6, 8).Value = 2
$ws.Cells.Item(6, 9).Value = 19.9
$ws.Cells.Item(6, 10).Value = 39.8

# Row 7
$ws.Cells.Item(7, 1).Value = "F35363_1556815384"
$ws.Cells.Item(7, 2).Value = 43587.73827284455
$ws.Cells.Item(7, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(7, 3).Value = "'123"
$ws.Cells.Item(7, 4).Value = "Carteiras"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = "'02"
$ws.Cells.Item(7, 7).Value = "Marca: Element; Tamanho: Unic; "
$ws.Cells.Item(7, 8).Value = 1
$ws.Cells.Item(7, 9).Value = 24.9
$ws.Cells.Item(7, 10).Value = 24.9

# Row 8
$ws.Cells.Item(8, 1).Value = "F42995_1556816158"
$ws.Cells.Item(8, 2).Value = 43587.74721260411
$ws.Cells.Item(8, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(8, 3).Value = "'123"
$ws.Cells.Item(8, 4).Value = "Bonés"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = "'02"
$ws.Cells.Item(8, 7).Value = "Nome: Boné New Era 940 Leag Basic; "
$ws.Cells.Item(8, 8).Value = 2
$ws.Cells.Item(8, 9).Value = 18.9
$ws.Cells.Item(8, 10).Value = 37.8
